# Fix test data inline with linter: replace "birth_date" with "birth date"
# in cell C1 (header row) of every worksheet that has that column.

$wb = $excel.ActiveWorkbook

$sheetNames = @("Animal", "NamedThing", "Person", "Animal1", "NamedThing1", "Person1")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Cells.Item(1, 3).Value = "birth date"
}
